# khaoSatPzem004.xlsx - "Modify khaosatPzem, co the lay data theo dia chi"
# Adds a new "dia chi 46" (address 46) row, renames the existing
# "Vao mode set parameter" read/reply pair to a generic "comand Read" /
# "set paramete" pair, and appends two new read/reply rows for address 46.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ngay1552019")

# --- New row 16: address info ---------------------------------------------
$ws.Range("A16").Value = "dia chi 46"
$ws.Range("B16").Value = "2e 04 00 00 00 0a 77 92"

# --- Row 22: rename the old "Vao mode set parameter" label -----------------
$ws.Range("B22").Value = " parameter - comand Read"
$ws.Range("D22").Value = "rep Vao mode set parameter"

# --- Row 23: unchanged request/reply bytes (kept as-is) --------------------
$ws.Range("B23").Value = "f8 03 00 00 00 07 10 61"
$ws.Range("D23").Value = "f8 03 0e 00 01 59 d8 00 01 25 80 08 98 27 10 55 f0 7b 04  "

# --- Row 25: new "set paramete" / "Add 46, AP 23000" pair ------------------
$ws.Range("B25").Value = "set paramete"
$ws.Range("D25").Value = "Add 46, AP 23000"

# --- Row 26: new read command / reply bytes for address 46 -----------------
$ws.Range("B26").Value = "f8 10 00 00 00 04 08 01 00 d8 59 2e 00 80 25 e9 12"
$ws.Range("D26").Value = "f8 10 00 00 00 04 d5 a3"

# --- Row 29: additional read command / reply bytes --------------------------
$ws.Range("B29").Value = "f8 10 00 00 00 04 08 01 00 d8 59 57 00 80 25 f0 4e"
$ws.Range("D29").Value = "f8 10 00 00 00 04 d5 a3                                "

# --- View state: scroll/select like the edited workbook --------------------
$ws.Activate()
$ws.Range("A10").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B28").Select() | Out-Null
